$wb = $excel.ActiveWorkbook

# ---- Sheet ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 51
$ws.Range("H51").Value = 4575
$ws.Range("I51").Value = 3266.6667
$ws.Range("J51").Value = 5360
$ws.Range("K51").Value = 3266.6667
$ws.Range("L51").Value = 5360
$ws.Range("M51").Value = -2782.6667
$ws.Range("N51").Value = -6328
# Row 113
$ws.Range("H113").Value = 3645.6365
$ws.Range("I113").Value = 2666.2856
$ws.Range("J113").Value = 4102.6665
$ws.Range("K113").Value = 2666.2856
$ws.Range("L113").Value = 4102.6665
$ws.Range("M113").Value = 587.7143999999998
$ws.Range("N113").Value = -10610.6665
# Row 116
$ws.Range("H116").Value = 3174.0322
$ws.Range("I116").Value = 2141.5908
$ws.Range("J116").Value = 5697.778
$ws.Range("K116").Value = 2141.5908
$ws.Range("L116").Value = 5697.778
$ws.Range("M116").Value = 1300.4092
$ws.Range("N116").Value = -12581.778
# Row 138
$ws.Range("H138").Value = 1940.3231
$ws.Range("I138").Value = 1778.4333
$ws.Range("J138").Value = 2079.0857
$ws.Range("K138").Value = 5335.2999
$ws.Range("L138").Value = 6237.257100000001
$ws.Range("M138").Value = -195.2999
$ws.Range("N138").Value = -16517.2571

# ---- Sheet ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 44
$ws.Range("H44").Value = 40000
$ws.Range("J44").Value = 40000
$ws.Range("L44").Value = 40000
$ws.Range("N44").Value = -40976
# Row 55
$ws.Range("H55").Value = 25000
$ws.Range("I55").Value = 0
$ws.Range("J55").Value = 25000
$ws.Range("K55").Value = 0
$ws.Range("M55").Value = 25000
$ws.Range("N55").Value = -25630
$ws.Range("L55").ClearContents()
# Row 61
$ws.Range("H61").Value = 2799.4614
$ws.Range("I61").Value = 1547.7059
$ws.Range("K61").Value = 1547.7059
$ws.Range("M61").Value = -1335.7059
# Row 63
$ws.Range("H63").Value = 3856.0908
$ws.Range("I63").Value = 1940.5555
$ws.Range("J63").Value = 12476
$ws.Range("K63").Value = 1940.5555
$ws.Range("L63").Value = 12476
$ws.Range("M63").Value = -1254.5555
$ws.Range("N63").Value = -13848
# Row 66
$ws.Range("H66").Value = 3856.0908
$ws.Range("I66").Value = 1940.5555
$ws.Range("J66").Value = 12476
$ws.Range("K66").Value = 9702.7775
$ws.Range("L66").Value = 62380
$ws.Range("M66").Value = -6270.7775
$ws.Range("N66").Value = -69244
# Row 102
$ws.Range("H102").Value = 1419
$ws.Range("I102").Value = 1468.5714
$ws.Range("J102").Value = 1245.5
$ws.Range("K102").Value = 1468.5714
$ws.Range("L102").Value = 1245.5
$ws.Range("M102").Value = 153.4286
$ws.Range("N102").Value = -4489.5
# Row 122
$ws.Range("H122").Value = 2488.1177
$ws.Range("I122").Value = 1253.2
$ws.Range("J122").Value = 11750
$ws.Range("K122").Value = 3759.6
$ws.Range("L122").Value = 35250
$ws.Range("M122").Value = -1309.6
$ws.Range("N122").Value = -40150
# Row 132
$ws.Range("H132").Value = 6587809.5
$ws.Range("I132").Value = 5740.8
$ws.Range("J132").Value = 8938548
$ws.Range("K132").Value = 17222.4
$ws.Range("L132").Value = 26815644
$ws.Range("M132").Value = -14692.4
$ws.Range("N132").Value = -26820704
# Row 136
$ws.Range("H136").Value = 2799.4614
$ws.Range("I136").Value = 1547.7059
$ws.Range("K136").Value = 4643.1177
$ws.Range("M136").Value = -2093.1177

# ---- Sheet BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 105
$ws.Range("H105").Value = 3308.75
$ws.Range("I105").Value = 1840
$ws.Range("J105").Value = 6540
$ws.Range("K105").Value = 1840
$ws.Range("L105").Value = 6540
$ws.Range("M105").Value = -93
$ws.Range("N105").Value = -10034
# Row 134
$ws.Range("H134").Value = 6396.8286
$ws.Range("I134").Value = 2918.1177
$ws.Range("J134").Value = 9682.277
$ws.Range("K134").Value = 8754.3531
$ws.Range("L134").Value = 29046.831
$ws.Range("M134").Value = -6219.3531
$ws.Range("N134").Value = -34116.831

# ---- Sheet CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 58
$ws.Range("H58").Value = 2789099.2
$ws.Range("I58").Value = 4687.3335
$ws.Range("J58").Value = 5573511
$ws.Range("K58").Value = 4687.3335
$ws.Range("L58").Value = 5573511
$ws.Range("M58").Value = -4484.3335
$ws.Range("N58").Value = -5573917
# Row 136
$ws.Range("H136").Value = 2789099.2
$ws.Range("I136").Value = 4687.3335
$ws.Range("J136").Value = 5573511
$ws.Range("K136").Value = 14062.0005
$ws.Range("L136").Value = 16720533
$ws.Range("M136").Value = -11512.0005
$ws.Range("N136").Value = -16725633
# Row 137
$ws.Range("H137").Value = 0
$ws.Range("J137").Value = 0
$ws.Range("N137").Value = 0
$ws.Range("L137").ClearContents()

# ---- Sheet GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 102
$ws.Range("H102").Value = 2465594.2
$ws.Range("I102").Value = 3573204.8
$ws.Range("J102").Value = 4237.3335
$ws.Range("K102").Value = 3573204.8
$ws.Range("L102").Value = 4237.3335
$ws.Range("M102").Value = -3571582.8
$ws.Range("N102").Value = -7481.3335
# Row 113
$ws.Range("H113").Value = 2777.1667
$ws.Range("I113").Value = 2783.3333
$ws.Range("J113").Value = 2771
$ws.Range("K113").Value = 2783.3333
$ws.Range("L113").Value = 2771
$ws.Range("M113").Value = -613.3332999999998
$ws.Range("N113").Value = -7111
# Row 116
$ws.Range("H116").Value = 0
$ws.Range("J116").Value = 0
$ws.Range("N116").Value = 0
$ws.Range("L116").ClearContents()
# Row 132
$ws.Range("H132").Value = 3363
$ws.Range("I132").Value = 5907
$ws.Range("J132").Value = 2599.8
$ws.Range("K132").Value = 17721
$ws.Range("L132").Value = 7799.400000000001
$ws.Range("M132").Value = -15191
$ws.Range("N132").Value = -12859.4

# ---- Sheet LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 22
$ws.Range("H22").Value = 1252.4706
$ws.Range("I22").Value = 999.1111
$ws.Range("J22").Value = 1537.5
$ws.Range("K22").Value = 999.1111
$ws.Range("L22").Value = 1537.5
$ws.Range("M22").Value = -704.1111
$ws.Range("N22").Value = -2127.5
# Row 27
$ws.Range("H27").Value = 1252.4706
$ws.Range("I27").Value = 999.1111
$ws.Range("J27").Value = 1537.5
$ws.Range("K27").Value = 999.1111
$ws.Range("L27").Value = 1537.5
$ws.Range("M27").Value = -892.1111
$ws.Range("N27").Value = -1751.5
# Row 40
$ws.Range("H40").Value = 101003480
$ws.Range("I40").Value = 168334500
$ws.Range("J40").Value = 6951.25
$ws.Range("K40").Value = 168334500
$ws.Range("L40").Value = 6951.25
$ws.Range("M40").Value = -168334364
$ws.Range("N40").Value = -7223.25
# Row 46
$ws.Range("H46").Value = 834440.25
$ws.Range("I46").Value = 832.3333
$ws.Range("J46").Value = 1668048.1
$ws.Range("K46").Value = 832.3333
$ws.Range("L46").Value = 1668048.1
$ws.Range("M46").Value = -644.3333
$ws.Range("N46").Value = -1668424.1
# Row 61
$ws.Range("H61").Value = 3550.7058
$ws.Range("I61").Value = 1576.2
$ws.Range("J61").Value = 6371.4287
$ws.Range("K61").Value = 1576.2
$ws.Range("L61").Value = 6371.4287
$ws.Range("M61").Value = -1374.2
$ws.Range("N61").Value = -6775.4287
# Row 113
$ws.Range("H113").Value = 3550.7058
$ws.Range("I113").Value = 1576.2
$ws.Range("J113").Value = 6371.4287
$ws.Range("K113").Value = 1576.2
$ws.Range("L113").Value = 6371.4287
$ws.Range("M113").Value = 593.8
$ws.Range("N113").Value = -10711.4287
# Row 122
$ws.Range("H122").Value = 19187.5
$ws.Range("J122").Value = 7250
$ws.Range("L122").Value = 21750
$ws.Range("N122").Value = -26650
# Row 136
$ws.Range("H136").Value = 45456452
$ws.Range("I136").Value = 83334700
$ws.Range("K136").Value = 250004100
$ws.Range("M136").Value = -250001550

# ---- Sheet WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 122
$ws.Range("H122").Value = 8268.695
$ws.Range("I122").Value = 9810
$ws.Range("J122").Value = 4745.7144
$ws.Range("K122").Value = 29430
$ws.Range("L122").Value = 14237.1432
$ws.Range("M122").Value = -26980
$ws.Range("N122").Value = -19137.1432
